$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = "329.06"
$ws.Range("E2").Value = "1.35%"
$rng.Style = "Normal"

$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = "41.22"
$ws.Range("E3").Value = "4.80%"
$rng.Style = "Normal"

$rng = $ws.Range("D4:E4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = "5.657"
$ws.Range("E4").Value = "-0.70%"
$rng.Style = "Normal"

$rng = $ws.Range("D5:E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = "0.08177"
$ws.Range("E5").Value = "2.15%"
$rng.Style = "Normal"

$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = "8.750"
$ws.Range("E6").Value = "1.82%"
$rng.Style = "Normal"

$rng = $ws.Range("D7:E7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = "2.005"
$ws.Range("E7").Value = "0.29%"
$rng.Style = "Normal"

$rng = $ws.Range("D8:E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = "4.483"
$ws.Range("E8").Value = "-0.25%"
$rng.Style = "Normal"

$rng = $ws.Range("D9:E9")
$rng.NumberFormat = "@"
$ws.Range("D9").Value = "2.945"
$ws.Range("E9").Value = "-0.96%"
$rng.Style = "Normal"

$rng = $ws.Range("D10:E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = "0.9206"
$ws.Range("E10").Value = "-0.56%"
$rng.Style = "Normal"

$rng = $ws.Range("D11:E11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = "0.1280"
$ws.Range("E11").Value = "3.45%"
$rng.Style = "Normal"

$rng = $ws.Range("D12:E12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = "0.1949"
$ws.Range("E12").Value = "-1.39%"
$rng.Style = "Normal"

$rng = $ws.Range("D13:E13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = "0.09267"
$ws.Range("E13").Value = "0.06%"
$rng.Style = "Normal"

$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = "0.03893"
$ws.Range("E14").Value = "7.88%"
$rng.Style = "Normal"

$rng = $ws.Range("E15")
$rng.NumberFormat = "@"
$ws.Range("E15").Value = "0.93%"
$rng.Style = "Normal"

$rng = $ws.Range("D16:E16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = "0.001299"
$ws.Range("E16").Value = "0.40%"
$rng.Style = "Normal"

$rng = $ws.Range("D17:E17")
$rng.NumberFormat = "@"
$ws.Range("D17").Value = "0.006224"
$ws.Range("E17").Value = "1.44%"
$rng.Style = "Normal"

$rng = $ws.Range("E19")
$rng.NumberFormat = "@"
$ws.Range("E19").Value = "2.94%"
$rng.Style = "Normal"

$rng = $ws.Range("D21:E21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = "8.244"
$ws.Range("E21").Value = "-5.44%"
$rng.Style = "Normal"

$rng = $ws.Range("E22")
$rng.NumberFormat = "@"
$ws.Range("E22").Value = "0.05%"
$rng.Style = "Normal"

$rng = $ws.Range("E23")
$rng.NumberFormat = "@"
$ws.Range("E23").Value = "-0.09%"
$rng.Style = "Normal"

$rng = $ws.Range("D24:E24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = "0.04403"
$ws.Range("E24").Value = "-0.08%"
$rng.Style = "Normal"

$rng = $ws.Range("D25:E25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = "0.001258"
$ws.Range("E25").Value = "-0.16%"
$rng.Style = "Normal"

$rng = $ws.Range("D26:E26")
$rng.NumberFormat = "@"
$ws.Range("D26").Value = "0.004307"
$ws.Range("E26").Value = "-6.97%"
$rng.Style = "Normal"

$rng = $ws.Range("D27:E27")
$rng.NumberFormat = "@"
$ws.Range("D27").Value = "0.0001202"
$ws.Range("E27").Value = "4.48%"
$rng.Style = "Normal"

$rng = $ws.Range("D39:E39")
$rng.NumberFormat = "@"
$ws.Range("D39").Value = "0.02792"
$ws.Range("E39").Value = "12.14%"
$rng.Style = "Normal"

$rng = $ws.Range("D40:E40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = "0.05397"
$ws.Range("E40").Value = "1.16%"
$rng.Style = "Normal"

$rng = $ws.Range("D41:E41")
$rng.NumberFormat = "@"
$ws.Range("D41").Value = "0.007792"
$ws.Range("E41").Value = "4.42%"
$rng.Style = "Normal"

$rng = $ws.Range("D42:E42")
$rng.NumberFormat = "@"
$ws.Range("D42").Value = "0.1416"
$ws.Range("E42").Value = "0.75%"
$rng.Style = "Normal"

$rng = $ws.Range("D43:E43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = "0.008940"
$ws.Range("E43").Value = "-7.02%"
$rng.Style = "Normal"

$rng = $ws.Range("D44:E44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = "0.002173"
$ws.Range("E44").Value = "2.69%"
$rng.Style = "Normal"

$rng = $ws.Range("D45:E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = "0.01149"
$ws.Range("E45").Value = "11.26%"
$rng.Style = "Normal"

$rng = $ws.Range("E46")
$rng.NumberFormat = "@"
$ws.Range("E46").Value = "-1.95%"
$rng.Style = "Normal"

$rng = $ws.Range("D47:E47")
$rng.NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "0.06%"
$rng.Style = "Normal"

$rng = $ws.Range("D48:E48")
$rng.NumberFormat = "@"
$ws.Range("D48").Value = "0.003212"
$ws.Range("E48").Value = "8.13%"
$rng.Style = "Normal"

$rng = $ws.Range("D50:E50")
$rng.NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").Value = "0.06%"
$rng.Style = "Normal"

$rng = $ws.Range("D51:E51")
$rng.NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").Value = "0.06%"
$rng.Style = "Normal"
